# "Generate Report for Handback"
#
# This script brings the localization-status workbook up to date after a
# handback run completed for the "a.md" source file:
#   - the status text for a.md moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview sheet + per-locale sheets)
#   - the per-locale sheets (zh-cn, de-de) get their "Latest Target File" and
#     "Latest Handback File" columns filled in, with a hyperlink on the
#     target file name, plus a real "Latest Handback DateTime" timestamp
#   - a couple of columns are widened so the new values are readable

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (shows up on the Overview summary sheet as well as each locale sheet)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Helper: apply the same blue-underline look used by the existing hyperlink
# cells (A2/A3) to a newly-populated hyperlink cell.
function Set-HyperlinkLook($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# Grab the real URL already used for the a.md hyperlink on each sheet so we
# don't have to hard-code the repository address.
function Get-AmdUrl($ws) {
    $url = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Row -eq 2 -and $h.Range.Column -eq 1) {
            $url = $h.Address
        }
    }
    return $url
}

# --- zh-cn sheet: fill in Latest Target File / Latest Handback File / DateTime ---
$zhCnUrl = Get-AmdUrl $wsZhCn

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $zhCnUrl, "", "", "a.md")
Set-HyperlinkLook $wsZhCn.Range("I2")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhCnUrl, "", "", "a.md")
Set-HyperlinkLook $wsZhCn.Range("I3")

$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-29 10:36:32"
$wsZhCn.Range("K3").Value = "2016-08-29 10:36:32"

# --- de-de sheet: same fields, different handback file/time ---
$deDeUrl = Get-AmdUrl $wsDeDe

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $deDeUrl, "", "", "a.md")
Set-HyperlinkLook $wsDeDe.Range("I2")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deDeUrl, "", "", "a.md")
Set-HyperlinkLook $wsDeDe.Range("I3")

$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-29 10:36:38"
$wsDeDe.Range("K3").Value = "2016-08-29 10:36:38"

# --- Column widths: widen the status column(s) and the Latest Handback File
#     column so the longer text fits ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.8   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.8   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth = 29.8    # C: Status
$wsZhCn.Columns.Item(10).ColumnWidth = 39.2   # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth = 29.8    # C: Status
$wsDeDe.Columns.Item(10).ColumnWidth = 39.2   # J: Latest Handback File

Write-Output "Handback report generated"
